# Exceloperations data files added
#
# - login sheet (Sheet1): E3 gets a "Aborted" result value.
# - Sheet2 becomes a "TestCase Name / UserName / Password" lookup table
#   (TC-1..TC-4 mapped to the Admin/User1/User2/User3 creds already on
#   the login sheet).
# - Selection / active-tab bookkeeping ends up on Sheet2.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---- login sheet: fill in the previously-empty Result cell ----
$ws1.Range("E3").Value = "Aborted"

# ---- Sheet2: new TestCase Name / UserName / Password table ----
$ws2.Range("A1").Value = "TestCase Name"
$ws2.Range("B1").Value = "UserName"
$ws2.Range("C1").Value = "Password"

$ws2.Range("A2").Value = "TC-1"
$ws2.Range("B2").Value = "Admin"
$ws2.Range("C2").Value = "adminpass"

$ws2.Range("A3").Value = "TC-2"
$ws2.Range("B3").Value = "User1"
$ws2.Range("C3").Value = "user1pass"

$ws2.Range("A4").Value = "TC-3"
$ws2.Range("B4").Value = "User2"
$ws2.Range("C4").Value = "user2pass"

$ws2.Range("A5").Value = "TC-4"
$ws2.Range("B5").Value = "User3"
$ws2.Range("C5").Value = "user3pass"

# Re-use the same header / data styles already defined on the login sheet
# (bold centered header for A1:B1, bold header for C1, bordered body cells
# for A2:C5) instead of minting brand-new style entries.
$ws1.Range("A1").Copy() | Out-Null
$ws2.Range("A1:B1").PasteSpecial(-4122) | Out-Null

$ws1.Range("C1").Copy() | Out-Null
$ws2.Range("C1").PasteSpecial(-4122) | Out-Null

$ws1.Range("A2").Copy() | Out-Null
$ws2.Range("A2:C5").PasteSpecial(-4122) | Out-Null

# Column widths for the new sheet
$ws2.Columns.Item(1).ColumnWidth = 19.666666666666668
$ws2.Columns.Item(2).ColumnWidth = 11.833333333333334
$ws2.Columns.Item(3).ColumnWidth = 20.666666666666668

# ---- selection / active sheet bookkeeping ----
$ws1.Activate()
$ws1.Rows.Item(2).Select() | Out-Null

$ws2.Activate()
$excel.ActiveWindow.Zoom = 150
$ws2.Range("A4").Select() | Out-Null
